$wb = $excel.ActiveWorkbook

# Rename the two "gathered" sheets to reflect the new MVC Path/PathController
# naming scheme used by the refactored test harness.
$wsGathered0 = $wb.Worksheets.Item("onglet 1")
$wsGathered0.Name = "tab_column_gathered_0"

$wsGathered1 = $wb.Worksheets.Item("onglet 2")
$wsGathered1.Name = "tab_column_gathered_1"
